$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76 (shifts rows 76..209 down to 77..210)
$ws.Rows(76).Insert()

# Populate the newly inserted row 76 with this week's new record
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = "Femacal de La Calera"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44469
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 100112040
$ws.Range("G76").Value = "Cilantro"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 280
$ws.Range("K76").Value = 2500
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = 2786
$ws.Range("N76").Value = "$/docena de atados (3 kilos)"
$ws.Range("O76").Value = "Provincia de Quillota"
$ws.Range("P76").Value = 929
$ws.Range("Q76").Value = 3
$ws.Range("R76").Value = "Hortaliza"
